$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section: "Output to ADC filter" - low-pass filter added on the outputs.

# Row 20 - section header
$ws.Range("B20").Value = "Output to ADC filter"

# Row 22 - Resistor value (reuse the "Comma" number format already used by D10)
$ws.Range("D10").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("C22").Value = "R"
$ws.Range("D22").Value = 4700

# Row 23 - Capacitor value (reuse the scientific number format already used by D11/D13)
$ws.Range("D11").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("C23").Value = "C"
$ws.Range("D23").Value = 0.000001

# Row 25 - Cutoff frequency formula (reuse the "Comma" number format again)
$ws.Range("D10").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("C25").Value = "fc"
$ws.Range("D25").Formula = "=1/(2*PI()*D22*D23)"

$excel.CutCopyMode = $false
